# Recursos.xlsx edit: rename sheets, fix duplicate "Contruccion"/"Construccion"
# shared string, and update Sheet1's active selection.

$wb = $excel.ActiveWorkbook

# --- Rename sheets ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Fix the misspelled duplicate shared string ("Contruccion") ------
# G10 was incorrectly pointing at the misspelled "Contruccion" entry;
# correcting its value to "Construccion" removes the now-unused
# duplicate string from the shared string table (uniqueCount 118 -> 117)
# and every subsequent shared-string index shifts down by one.
$ws1.Range("G10").Value = "Construccion"

$ws1.Name = "Recursos Naturales"
$ws2.Name = "Productos"

# --- Update the active selection on Sheet1 ----------------------------
$ws1.Select() | Out-Null
$ws1.Range("Q4").Select() | Out-Null
